# Generate Report for Handoff
#
# The handoff XLIFF for "1437a34b-8d45-4312-90bf-9e34f9defe7c.md" was
# (re)generated, bumping its "Latest Handoff" timestamps on the Overview
# sheet and on each per-locale (zh-cn / de-de) status sheet.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for the row
# belonging to 1437a34b-8d45-4312-90bf-9e34f9defe7c.md (row 5).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G5").Value = "2016-08-22 16:44:21"

# zh-cn sheet: "Latest Handoff Datetime" column (H) for the same file (row 5).
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H5").Value = "2016-08-22 16:44:17"

# de-de sheet: "Latest Handoff Datetime" column (H) for the same file (row 5).
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H5").Value = "2016-08-22 16:44:21"
